$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Long Selenium "no such element" error text shared by F11:F14 (new shared string)
$longError = @'
no such element: Unable to locate element: {"method":"css selector","selector":"#lblServiceID"}
  (Session info: chrome-headless-shell=121.0.6167.189)
For documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.130.41', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome-headless-shell, browserVersion: 121.0.6167.189, chrome: {chromedriverVersion: 121.0.6167.184 (057a8ae7deb..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:50375}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 4e0b85e0b26caa43044a6b7ed4319fa7
*** Element info: {Using=id, value=lblServiceID}
'@

function Set-TextValue {
    param($range, [string]$text)
    # Assign as text (leading apostrophe forces Excel to treat a
    # numeric-looking string as text instead of a number) and strip any
    # auto-applied number formatting so the cell keeps using the
    # workbook's default style (important when the text looks like a
    # plain number, e.g. service IDs).
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Column C - Service ID values
Set-TextValue $ws.Range("C2")  "10276024"
Set-TextValue $ws.Range("C3")  "10276026"
Set-TextValue $ws.Range("C4")  "10276032"
Set-TextValue $ws.Range("C5")  "10276039"
Set-TextValue $ws.Range("C11") "10276088"
Set-TextValue $ws.Range("C12") "10276068"
Set-TextValue $ws.Range("C13") "10276077"
Set-TextValue $ws.Range("C14") "10276082"
Set-TextValue $ws.Range("C24") "137010051"

# Column F - Result / error message values
$ws.Range("F2").Value = 'Cannot invoke "org.openqa.selenium.WebElement.isDisplayed()" because "element" is null'
$ws.Range("F3").Value = 'Cannot invoke "org.openqa.selenium.WebElement.getText()" because the return value of "connect_OCBaseMethods.TCAcknowledge.isElementPresent(String)" is null'

$ws.Range("F11").Value = $longError
$ws.Range("F12").Value = $longError
$ws.Range("F13").Value = $longError
$ws.Range("F14").Value = $longError
